# Auto-generated edit script: updates market-price derived columns (H-N)
# on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 70300.39
$ws.Range("I76").Value = 103778.95
$ws.Range("J76").Value = 3343.2727
$ws.Range("K76").Value = 103778.95
$ws.Range("L76").Value = 3343.2727
$ws.Range("M76").Value = -103463.95
$ws.Range("N76").Value = -3973.2727

$ws.Range("H79").Value = 70300.39
$ws.Range("I79").Value = 103778.95
$ws.Range("J79").Value = 3343.2727
$ws.Range("K79").Value = 103778.95
$ws.Range("L79").Value = 3343.2727
$ws.Range("M79").Value = -102686.95
$ws.Range("N79").Value = -5527.2727

$ws.Range("H80").Value = 1083.9722
$ws.Range("I80").Value = 516.6539
$ws.Range("J80").Value = 2559
$ws.Range("K80").Value = 1549.9617
$ws.Range("L80").Value = 7677
$ws.Range("M80").Value = -551.9617000000001
$ws.Range("N80").Value = -9673

$ws.Range("H83").Value = 1083.9722
$ws.Range("I83").Value = 516.6539
$ws.Range("J83").Value = 2559
$ws.Range("K83").Value = 4649.8851
$ws.Range("L83").Value = 23031
$ws.Range("M83").Value = 342.1148999999996
$ws.Range("N83").Value = -33015

$ws.Range("H86").Value = 84394.78
$ws.Range("I86").Value = 125665.5
$ws.Range("J86").Value = 1853.3334
$ws.Range("K86").Value = 125665.5
$ws.Range("L86").Value = 1853.3334
$ws.Range("M86").Value = -124542.5
$ws.Range("N86").Value = -4099.3334

$ws.Range("H89").Value = 84394.78
$ws.Range("I89").Value = 125665.5
$ws.Range("J89").Value = 1853.3334
$ws.Range("K89").Value = 628327.5
$ws.Range("L89").Value = 9266.666999999999
$ws.Range("M89").Value = -622711.5
$ws.Range("N89").Value = -20498.667

$ws.Range("H92").Value = 617.8095
$ws.Range("I92").Value = 594.3077
$ws.Range("J92").Value = 656
$ws.Range("K92").Value = 594.3077
$ws.Range("L92").Value = 656
$ws.Range("M92").Value = 653.6923
$ws.Range("N92").Value = -3152

$ws.Range("H98").Value = 1698.3334
$ws.Range("I98").Value = 1638
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 1638
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -140
$ws.Range("N98").Value = -4996

$ws.Range("H122").Value = 1698.3334
$ws.Range("I122").Value = 1638
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4914
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2464
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5775.46
$ws.Range("I32").Value = 4617.859
$ws.Range("J32").Value = 12335.2
$ws.Range("K32").Value = 4617.859
$ws.Range("L32").Value = 12335.2
$ws.Range("M32").Value = -4330.859
$ws.Range("N32").Value = -12909.2

$ws.Range("H74").Value = 187090.12
$ws.Range("I74").Value = 1853.159
$ws.Range("J74").Value = 1002132.8
$ws.Range("K74").Value = 1853.159
$ws.Range("L74").Value = 1002132.8
$ws.Range("M74").Value = -979.1590000000001
$ws.Range("N74").Value = -1003880.8

$ws.Range("H77").Value = 187090.12
$ws.Range("I77").Value = 1853.159
$ws.Range("J77").Value = 1002132.8
$ws.Range("K77").Value = 9265.795
$ws.Range("L77").Value = 5010664
$ws.Range("M77").Value = -4897.795
$ws.Range("N77").Value = -5019400

$ws.Range("H132").Value = 23694.89
$ws.Range("I132").Value = 30082.742
$ws.Range("K132").Value = 90248.226
$ws.Range("M132").Value = -87718.226

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1582.5
$ws.Range("J20").Value = 1866.6666
$ws.Range("L20").Value = 1866.6666
$ws.Range("N20").Value = -2360.6666

$ws.Range("H86").Value = 6898629.5
$ws.Range("I86").Value = 9092908
$ws.Range("J86").Value = 2327
$ws.Range("K86").Value = 9092908
$ws.Range("L86").Value = 2327
$ws.Range("M86").Value = -9091785
$ws.Range("N86").Value = -4573

$ws.Range("H89").Value = 6898629.5
$ws.Range("I89").Value = 9092908
$ws.Range("J89").Value = 2327
$ws.Range("K89").Value = 45464540
$ws.Range("L89").Value = 11635
$ws.Range("M89").Value = -45458924
$ws.Range("N89").Value = -22867

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2029.7542
$ws.Range("I31").Value = 1287.1852
$ws.Range("J31").Value = 2619.4412
$ws.Range("K31").Value = 1287.1852
$ws.Range("L31").Value = 2619.4412
$ws.Range("M31").Value = -992.1851999999999
$ws.Range("N31").Value = -3209.4412

$ws.Range("H34").Value = 2029.7542
$ws.Range("I34").Value = 1287.1852
$ws.Range("J34").Value = 2619.4412
$ws.Range("K34").Value = 1287.1852
$ws.Range("L34").Value = 2619.4412
$ws.Range("M34").Value = -1085.1852
$ws.Range("N34").Value = -3023.4412

$ws.Range("H86").Value = 279944.62
$ws.Range("I86").Value = 359135.94
$ws.Range("K86").Value = 359135.94
$ws.Range("M86").Value = -358012.94

$ws.Range("H89").Value = 279944.62
$ws.Range("I89").Value = 359135.94
$ws.Range("K89").Value = 1795679.7
$ws.Range("M89").Value = -1790063.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 662.9828
$ws.Range("I5").Value = 379.0645
$ws.Range("K5").Value = 1137.1935
$ws.Range("M5").Value = -1025.1935

$ws.Range("H122").Value = 747.4722
$ws.Range("I122").Value = 400.27274
$ws.Range("J122").Value = 1293.0714
$ws.Range("K122").Value = 3602.45466
$ws.Range("L122").Value = 11637.6426
$ws.Range("M122").Value = -1152.45466
$ws.Range("N122").Value = -16537.6426

$ws.Range("H123").Value = 2663.544
$ws.Range("I123").Value = 1137.1428
$ws.Range("J123").Value = 2877.24
$ws.Range("K123").Value = 3411.4284
$ws.Range("L123").Value = 8631.719999999999
$ws.Range("M123").Value = -961.4284000000002
$ws.Range("N123").Value = -13531.72

$ws.Range("H127").Value = 1202.1111
$ws.Range("J127").Value = 1202.1111
$ws.Range("L127").Value = 3606.3333
$ws.Range("N127").Value = -13526.3333

$ws.Range("H135").Value = 662.9828
$ws.Range("I135").Value = 379.0645
$ws.Range("K135").Value = 3411.5805
$ws.Range("M135").Value = -876.5805

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 51842240
$ws.Range("I70").Value = 207353440
$ws.Range("J70").Value = 5172.9165
$ws.Range("K70").Value = 207353440
$ws.Range("L70").Value = 5172.9165
$ws.Range("M70").Value = -207353170
$ws.Range("N70").Value = -5712.9165

$ws.Range("H73").Value = 51842240
$ws.Range("I73").Value = 207353440
$ws.Range("J73").Value = 5172.9165
$ws.Range("K73").Value = 207353440
$ws.Range("L73").Value = 5172.9165
$ws.Range("M73").Value = -207352504
$ws.Range("N73").Value = -7044.9165

$ws.Range("H122").Value = 1750.7273
$ws.Range("I122").Value = 2009.25
$ws.Range("J122").Value = 1440.5
$ws.Range("K122").Value = 6027.75
$ws.Range("L122").Value = 4321.5
$ws.Range("M122").Value = -3577.75
$ws.Range("N122").Value = -9221.5

$ws.Range("H132").Value = 2616.8
$ws.Range("I132").Value = 2052.8667
$ws.Range("K132").Value = 6158.6001
$ws.Range("M132").Value = -3628.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 412.0625
$ws.Range("J16").Value = 728
$ws.Range("L16").Value = 728
$ws.Range("N16").Value = -1068

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 21250.275
$ws.Range("J100").Value = 841.3158
$ws.Range("L100").Value = 1682.6316
$ws.Range("N100").Value = -2764.6316
